$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Collect all cells that need to stay as literal text (not auto-converted to number/percentage)
$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","E7","D8","E8","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","E18","D19","E19","D20","E20","E21","D22","E22","D23","E23","E24","D25","E25","E26","D27","E27","D39","E39","D40","E40","E41","D42","E42","D43","E43","E44","D45","E45","D46","E46","E47","D48","E48","D49","E49","D50","E50","E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '326.41'
$ws.Range("E2").Value = '-1.14%'
$ws.Range("D3").Value = '45.08'
$ws.Range("E3").Value = '3.05%'
$ws.Range("D4").Value = '5.571'
$ws.Range("E4").Value = '-5.81%'
$ws.Range("D5").Value = '0.08088'
$ws.Range("E5").Value = '-2.49%'
$ws.Range("D6").Value = '8.714'
$ws.Range("E6").Value = '-0.87%'
$ws.Range("E7").Value = '-3.64%'
$ws.Range("D8").Value = '1.910'
$ws.Range("E8").Value = '-2.25%'
$ws.Range("D10").Value = '0.9486'
$ws.Range("E10").Value = '1.90%'
$ws.Range("D11").Value = '0.1180'
$ws.Range("E11").Value = '-5.92%'
$ws.Range("D12").Value = '0.1895'
$ws.Range("E12").Value = '-2.34%'
$ws.Range("D13").Value = '0.1007'
$ws.Range("E13").Value = '6.15%'
$ws.Range("D14").Value = '0.04182'
$ws.Range("E14").Value = '5.33%'
$ws.Range("D15").Value = '0.1064'
$ws.Range("E15").Value = '0.13%'
$ws.Range("D16").Value = '0.001271'
$ws.Range("E16").Value = '-3.18%'
$ws.Range("D17").Value = '0.005985'
$ws.Range("E17").Value = '-0.26%'
$ws.Range("E18").Value = '2.40%'
$ws.Range("D19").Value = '0.3487'
$ws.Range("E19").Value = '-0.65%'
$ws.Range("D20").Value = '8.543'
$ws.Range("E20").Value = '-6.79%'
$ws.Range("E21").Value = '0.15%'
$ws.Range("D22").Value = '0.2664'
$ws.Range("E22").Value = '2.45%'
$ws.Range("D23").Value = '0.04275'
$ws.Range("E23").Value = '-3.36%'
$ws.Range("E24").Value = '-1.55%'
$ws.Range("D25").Value = '0.004599'
$ws.Range("E25").Value = '3.94%'
$ws.Range("E26").Value = '3.40%'
$ws.Range("D27").Value = '0.0004001'
$ws.Range("E27").Value = '0.14%'
$ws.Range("D39").Value = '0.02671'
$ws.Range("E39").Value = '-5.56%'
$ws.Range("D40").Value = '0.05555'
$ws.Range("E40").Value = '-1.59%'
$ws.Range("E41").Value = '24.80%'
$ws.Range("D42").Value = '0.007673'
$ws.Range("E42").Value = '-3.20%'
$ws.Range("D43").Value = '0.1394'
$ws.Range("E43").Value = '-1.96%'
$ws.Range("E44").Value = '-1.43%'
$ws.Range("D45").Value = '0.008706'
$ws.Range("E45").Value = '-0.60%'
$ws.Range("D46").Value = '0.00007128'
$ws.Range("E46").Value = '-2.37%'
$ws.Range("E47").Value = '0.14%'
$ws.Range("D48").Value = '0.003445'
$ws.Range("E48").Value = '-6.11%'
$ws.Range("D49").Value = '0.002277'
$ws.Range("E49").Value = '-0.19%'
$ws.Range("D50").Value = '0.00002106'
$ws.Range("E50").Value = '0.14%'
$ws.Range("E51").Value = '0.14%'
